$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to keep text formatting so numeric-looking price strings
# (e.g. "594.30", "0.999") are not coerced into numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.325.05"
$ws.Range("E2").Value = "  +4.19%  "
$ws.Range("D3").Value = "3.631.78"
$ws.Range("E3").Value = "  +3.93%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "594.30"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "195.63"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").Value = "3.625.47"
$ws.Range("E8").Value = "  +3.90%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +2.98%  "
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("D12").Value = "58.86"
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").Value = "4.212.81"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").Value = "19.84"
$ws.Range("E16").Value = "  +4.49%  "
$ws.Range("D17").Value = "3.638.12"
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").Value = "70.316.55"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("E21").Value = "  +4.21%  "
$ws.Range("D22").Value = "488.20"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +15.33%  "
$ws.Range("D24").Value = "5.34"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +7.00%  "
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").Value = "9.63"
$ws.Range("E29").Value = "  +4.75%  "
$ws.Range("D30").Value = "7.97"
$ws.Range("E30").Value = "  +10.25%  "
$ws.Range("D31").Value = "32.98"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("E32").Value = "  +8.08%  "
$ws.Range("D33").Value = "627.08"
$ws.Range("E33").Value = "  +5.08%  "
$ws.Range("D34").Value = "12.29"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("D35").Value = "65.96"
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("D36").Value = "40.70"
$ws.Range("E36").Value = "  +10.67%  "
$ws.Range("E37").Value = "  +6.79%  "
$ws.Range("D38").Value = "0.0₃0823"
$ws.Range("E38").Value = "  +6.66%  "
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "3.297.60"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  +6.99%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +10.84%  "
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").Value = "9.22"
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("D50").Value = "3.35"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("E51").Value = "  -0.16%  "
